## ---------------------------------------------------------------------
## Commit: "Tue, May 12, 2020  8:05:27 PM"
##
## Two logically separate edits:
##
## 1) Three tables (on the slides whose table grid widths are
##    2879725 / 3424250 / 2881325 EMU -- i.e. slides 14, 15 and 16)
##    get their table style switched from the deck's custom
##    "Table_0" style ({732938CB-28EC-4491-9528-9433618EBCEF}) to the
##    built-in PowerPoint "No Style, Table Grid" style
##    ({873121FD-0E87-461C-AD5D-959C63CB78D2}).
##
## 2) The presentation's Design ("Integral" / "Red Violet" colour
##    scheme) is swapped out for the stock "Office Theme" design --
##    i.e. every themed colour (dk1/lt1/dk2/lt2/accent1-6/hlink/
##    folHlink) used by the slide master's theme is reset back to the
##    Office defaults.
## ---------------------------------------------------------------------

$p = $ppt.ActivePresentation

## --- 1) retarget the three tables onto the built-in table style -------
$newTableStyleId = "{873121FD-0E87-461C-AD5D-959C63CB78D2}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

## --- 2) swap the design/theme colours back to the default Office theme
## Office Theme colour scheme (index -> RGB long, BGR-packed the way the
## PowerPoint object model stores Long colour values):
##   1 dk1       000000
##   2 lt1       FFFFFF
##   3 dk2       44546A
##   4 lt2       E7E6E6
##   5 accent1   5B9BD5
##   6 accent2   ED7D31
##   7 accent3   A5A5A5
##   8 accent4   FFC000
##   9 accent5   4472C4
##  10 accent6   70AD47
##  11 hlink     0563C1
##  12 folHlink  954F72
$officeThemeColors = @{
    1  = 0x000000
    2  = 0xFFFFFF
    3  = 0x6A5444
    4  = 0xE6E6E7
    5  = 0xD59B5B
    6  = 0x317DED
    7  = 0xA5A5A5
    8  = 0x00C0FF
    9  = 0xC47244
    10 = 0x47AD70
    11 = 0xC16305
    12 = 0x724F95
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in 1..12) {
    $themeColors.Item($idx).RGB = $officeThemeColors[$idx]
}
